$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price values stay as text, matching the source format
$ws.Range("D2").Value = "58.113.99"
$ws.Range("E2").Value = "  +0.44%  "
$ws.Range("D3").Value = "2.462.20"
$ws.Range("E3").Value = "  +0.42%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "510.67"
$ws.Range("E5").Value = "  -2.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.90"
$ws.Range("E6").Value = "  +3.16%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  -0.60%  "
$ws.Range("D9").Value = "2.461.03"
$ws.Range("E9").Value = "  +0.15%  "
$ws.Range("E10").Value = "  +0.60%  "
$ws.Range("E11").Value = "  -0.69%  "
$ws.Range("E12").Value = "  +0.81%  "
$ws.Range("E13").Value = "  -6.20%  "
$ws.Range("D14").Value = "2.895.55"
$ws.Range("E14").Value = "  +0.31%  "
$ws.Range("D15").Value = "57.971.29"
$ws.Range("E15").Value = "  +0.31%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.98"
$ws.Range("E16").Value = "  +1.94%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000135"
$ws.Range("E17").Value = "  +1.59%  "
$ws.Range("D18").Value = "2.475.24"
$ws.Range("E18").Value = "  +0.88%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.38"
$ws.Range("E20").Value = "  -0.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "314.93"
$ws.Range("E21").Value = "  +0.64%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.45"
$ws.Range("E22").Value = "  +5.52%  "
$ws.Range("E23").Value = "  -0.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.73"
$ws.Range("E24").Value = "  -2.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.45"
$ws.Range("E25").Value = "  +0.80%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.994"
$ws.Range("E26").Value = "  -0.66%  "
$ws.Range("E27").Value = "  -0.36%  "
$ws.Range("E28").Value = "  -4.98%  "
$ws.Range("E29").Value = "  +5.12%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "172.27"
$ws.Range("E30").Value = "  -1.25%  "
$ws.Range("E31").Value = "  +0.27%  "
$ws.Range("E32").Value = "  +0.14%  "
$ws.Range("E33").Value = "  +0.34%  "
$ws.Range("E34").Value = "  +1.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  +0.12%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.11"
$ws.Range("E37").Value = "  +1.62%  "
$ws.Range("E38").Value = "  +5.41%  "
$ws.Range("E39").Value = "  +3.50%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.76"
$ws.Range("E40").Value = "  +1.30%  "
$ws.Range("E41").Value = "  +1.98%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.810"
$ws.Range("E42").Value = "  +0.51%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "137.16"
$ws.Range("E43").Value = "  +11.20%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.41"
$ws.Range("E44").Value = "  +0.72%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.92"
$ws.Range("E45").Value = "  +2.78%  "
$ws.Range("E46").Value = "  -0.83%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "255.82"
$ws.Range("E47").Value = "  -0.50%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0918"
$ws.Range("E48").Value = "  -0.07%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0495"
$ws.Range("E49").Value = "  +0.71%  "
$ws.Range("E50").Value = "  +2.12%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.28"
$ws.Range("E51").Value = "  +1.71%  "
